$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 87
$ws.Range("I2").Value = 218
$ws.Range("J2").Value = 846
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 259
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 159
$ws.Range("P2").Value = 3
$ws.Range("R2").Value = 17
$ws.Range("S2").Value = 104
$ws.Range("T2").Value = 145
$ws.Range("V2").Value = 1416
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1390
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 20
$ws.Range("AA2").Value = 8
